# Fix Training Data Issue (#48)
# The BF column holds a per-row "Date" label that was mistakenly set to the
# literal season-folder string (e.g. "6-12-2007-08"); it should hold the
# actual ISO game date "2008-06-12" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "6-12-2007-08"
$newValue = "2008-06-12"

# BF2:BF31 hold the old value (BF1 is the "Date" header and is untouched).
$firstRow = 2
$lastRow = 31
$col = "BF"

$rng = $ws.Range($col + $firstRow + ":" + $col + $lastRow)

# Pre-format the range as Text so the ISO-looking "2008-06-12" string isn't
# silently reinterpreted as a date serial number on assignment.
$rng.NumberFormat = "@"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Range($col + $r)
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}

# Restore the default cell style so no stray number-format style lingers on
# these cells (matches the original file, which left them unstyled).
$rng.Style = "Normal"
